# Add a "Type" column to the Compartments sheet (inserted before the
# existing "Comments" column), matching the commit
# "add type attribute column for Compartment".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compartments")

# Insert a new column E ("Type"), shifting Comments/References right.
$ws.Columns.Item(5).Insert()
$ws.Cells.Item(1, 5).Value2 = "Type"

# Re-apply the AutoFilter so its range covers the new column (A1:G2).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:G2").AutoFilter()

# Update the workbook-level _FilterDatabase defined name for this sheet.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Compartments!_FilterDatabase") {
        $n.RefersTo = "=Compartments!`$A`$1:`$G`$2"
    }
}

# Make Compartments the active sheet/tab, with E2 selected, matching the
# new workbookView activeTab and sheetView selection in the diff.
[void]$ws.Activate()
[void]$ws.Range("E2").Select()
